$d = $word.ActiveDocument

# --- Paragraph 1: "June 6, 2019" -> "October 31, 2019" ---
# The original paragraph has four runs: "June" | " " | "6" | ", 2019".
# A plain Find/Replace (or Range.Text=) on this engine coalesces the
# edited run together with its neighbours, so we drop temporary
# "barrier" bookmarks at the run boundaries we want to keep intact;
# bookmarks stop the auto-merge, and we delete the temporary ones
# afterwards (this does not undo the split). The real "_GoBack"
# bookmark is (re)created exactly where Word leaves it after this
# edit: right after the newly typed "October".

$goBack = $d.Bookmarks.Add("_GoBack", $d.Range(4, 4))
$barrierA = $d.Bookmarks.Add("ZZ_barrierA", $d.Range(5, 5))
$barrierB = $d.Bookmarks.Add("ZZ_barrierB", $d.Range(6, 6))

$d.Range(0, 4).Text = "October"
$d.Range(8, 9).Text = "31"

$d.Bookmarks("ZZ_barrierA").Delete()
$d.Bookmarks("ZZ_barrierB").Delete()

# --- Paragraph with "Dear Editor, Scientific Reports" -> "...mBio" ---
$d.Content.Find.Execute("Scientific Reports", $true, $false, $false, $false, $false,
                         $true, 1, $false, "mBio", 2)
